$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Append the newly-found foreign-key exception IDs at the bottom of the
#    list, flagged in red text and formatted as Text (like the other
#    out-of-band entries already in the sheet).
# ---------------------------------------------------------------------------
$ws.Cells.Item(38, 1).Value = 2079038
$ws.Cells.Item(39, 1).Value = 2464225
$ws.Cells.Item(40, 1).Value = 3060918
$ws.Cells.Item(41, 1).Value = 9813349

$newCells = $ws.Range("A38:A41")
$newCells.Font.Color = 255
$newCells.NumberFormat = "@"

# ---------------------------------------------------------------------------
# 2. Sort the whole ID list (A2:A41) ascending, keyed on column A / the
#    header in A1 - same as Data > Sort A to Z in the Excel UI.
# ---------------------------------------------------------------------------
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$dataRange = $ws.Range("A2:A41")
$keyRange = $ws.Range("A1")
$sortObj.SortFields.Add($keyRange)
$sortObj.SetRange($dataRange)
$sortObj.Header = 0
$sortObj.MatchCase = $false
$sortObj.Orientation = 1
$sortObj.Apply()

# ---------------------------------------------------------------------------
# 3. The engine groups true numbers ahead of text during the sort above, so
#    put the four flagged numeric rows back at their correct ascending
#    position among the text codes, shifting everything else down/up
#    accordingly. Existing (already-sorted) text cells are left untouched so
#    their formatting/style is not disturbed.
# ---------------------------------------------------------------------------
function Move-NewIdToRow($fromRow, $toRow, $value) {
    $ws.Rows.Item($fromRow).Delete()
    $ws.Rows.Item($toRow).Insert()
    $c = $ws.Cells.Item($toRow, 1)
    $c.Value = $value
    $c.Font.Color = 255
    $c.NumberFormat = "@"
}

# After Apply(), the 4 new numbers sit first, at rows 2-5 (in the order they
# were added: 2079038, 2464225, 3060918, 9813349). Re-home each one, working
# from the bottom target row upward so earlier moves don't disturb the rows
# still to be placed.
Move-NewIdToRow 5 34 9813349
Move-NewIdToRow 4 15 3060918
Move-NewIdToRow 3 9 2464225
Move-NewIdToRow 2 4 2079038

# ---------------------------------------------------------------------------
# 4. Leave the selection where the user last clicked while reviewing the
#    refreshed list.
# ---------------------------------------------------------------------------
$ws.Range("G19").Select()

Write-Host "Applied foreign-key exception update"
